$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.994.19"
$ws.Range("E2").Value = "  -0.85%  "

$ws.Range("D3").Value = "3.522.56"
$ws.Range("E3").Value = "  +0.70%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.03"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.99"
$ws.Range("E6").Value = "  -2.47%  "

$ws.Range("D7").Value = "3.520.65"
$ws.Range("E7").Value = "  +0.71%  "

$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.479"
$ws.Range("E9").Value = "  -1.27%  "

$ws.Range("E10").Value = "  -0.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.87"
$ws.Range("E11").Value = "  +3.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.423"
$ws.Range("E12").Value = "  -1.82%  "

$ws.Range("E13").Value = "  -0.39%  "

$ws.Range("D14").Value = "4.122.97"
$ws.Range("E14").Value = "  +0.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.68"
$ws.Range("E15").Value = "  -2.15%  "

$ws.Range("D16").Value = "3.517.19"
$ws.Range("E16").Value = "  +0.62%  "

$ws.Range("D17").Value = "67.022.59"
$ws.Range("E17").Value = "  -0.87%  "

$ws.Range("E18").Value = "  -0.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.84"
$ws.Range("E19").Value = "  +10.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.39"
$ws.Range("E20").Value = "  -2.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.36"
$ws.Range("E21").Value = "  -1.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "436.60"
$ws.Range("E22").Value = "  -2.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.611"
$ws.Range("E23").Value = "  -3.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.78"
$ws.Range("E24").Value = "  +2.23%  "

$ws.Range("D25").Value = "3.663.48"
$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000121"
$ws.Range("E27").Value = "  -3.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.86"
$ws.Range("E28").Value = "  -1.92%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.42"
$ws.Range("E29").Value = "  -3.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.50"
$ws.Range("E30").Value = "  +0.16%  "

$ws.Range("E31").Value = "  -2.76%  "

$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.20%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.42"
$ws.Range("E34").Value = "  -0.97%  "

$ws.Range("D35").Value = "3.516.84"
$ws.Range("E35").Value = "  +0.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.80"
$ws.Range("E36").Value = "  -2.85%  "

$ws.Range("E37").Value = "  -4.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.05"
$ws.Range("E38").Value = "  +0.72%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.08%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0893"
$ws.Range("E41").Value = "  +0.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "170.48"
$ws.Range("E42").Value = "  -2.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.44"
$ws.Range("E43").Value = "  -0.11%  "

$ws.Range("E44").Value = "  -10.17%  "

$ws.Range("E45").Value = "  +0.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.93"
$ws.Range("E46").Value = "  -0.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.33"
$ws.Range("E47").Value = "  +2.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "28.18"
$ws.Range("E48").Value = "  -6.13%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.48"
$ws.Range("E49").Value = "  -1.73%  "

$ws.Range("E50").Value = "  -2.63%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.994"
$ws.Range("E51").Value = "  +0.19%  "
